$wb = $excel.ActiveWorkbook

# Add a new worksheet after the last existing sheet, and name it "test_data2"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "test_data2"

# Populate the simple key/value table starting at C5
$data = @(
    @("key1", "value1"),
    @("key2", "value2"),
    @("key3", "value3"),
    @("key4", "value4")
)

$row = 5
foreach ($pair in $data) {
    $ws.Cells.Item($row, 3).Value = $pair[0]
    $ws.Cells.Item($row, 4).Value = $pair[1]
    $row++
}
